$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "ABHISHEK"
$ws.Range("B5").Value = "ABHISHEK"
$ws.Range("C5").Value = $true

$ws.Range("C5").Select() | Out-Null
